$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date/time style from A2 onto the new rows (A8:A25) before setting values
$ws.Range("A2").Copy()
$ws.Range("A8:A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate Time (A), Upper3 (B) / MA (C) values row by row
$ws.Range("A2").Value = 44125.625
$ws.Range("B2").Value = 2.107675754536493
$ws.Range("A3").Value = 44125.75
$ws.Range("C3").Value = 1.948202576129336
$ws.Range("A4").Value = 44142.07291666666
$ws.Range("B4").Value = 2.250337282945386
$ws.Range("A5").Value = 44142.58333333334
$ws.Range("C5").Value = 2.245845942664026
$ws.Range("A6").Value = 44178.51041666666
$ws.Range("B6").Value = 1.607061071303064
$ws.Range("A7").Value = 44178.83333333334
$ws.Range("C7").Value = 1.572894510843002
$ws.Range("A8").Value = 44180.125
$ws.Range("B8").Value = 1.719168800896331
$ws.Range("A9").Value = 44180.59375
$ws.Range("C9").Value = 1.718387281079325
$ws.Range("A10").Value = 44204.41666666666
$ws.Range("B10").Value = 3.120130548023874
$ws.Range("A11").Value = 44204.83333333334
$ws.Range("C11").Value = 3.162118492340033
$ws.Range("A12").Value = 44229.35416666666
$ws.Range("B12").Value = 5.351561416378089
$ws.Range("A13").Value = 44229.64583333334
$ws.Range("C13").Value = 5.071733591711183
$ws.Range("A14").Value = 44246.51041666666
$ws.Range("B14").Value = 10.35509668954349
$ws.Range("A15").Value = 44247.03125
$ws.Range("C15").Value = 11.08366481300376
$ws.Range("A16").Value = 44249.8125
$ws.Range("B16").Value = 13.34829008640281
$ws.Range("A17").Value = 44250.20833333334
$ws.Range("C17").Value = 13.26324218499715
$ws.Range("A18").Value = 44266.72916666666
$ws.Range("B18").Value = 15.67519183306398
$ws.Range("A19").Value = 44267.3125
$ws.Range("C19").Value = 16.01936378958497
$ws.Range("A20").Value = 44277.17708333334
$ws.Range("B20").Value = 15.45396416134136
$ws.Range("A21").Value = 44277.5625
$ws.Range("C21").Value = 15.76496271752683
$ws.Range("A22").Value = 44403.0625
$ws.Range("B22").Value = 30.71234535161814
$ws.Range("A23").Value = 44403.33333333334
$ws.Range("C23").Value = 29.36151843268062
$ws.Range("A24").Value = 44423.73958333334
$ws.Range("B24").Value = 52.47724936166023
$ws.Range("A25").Value = 44424.82291666666
$ws.Range("C25").Value = 64.31932694245472
